$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell holding the default (unstyled) style, used to restore style
# after temporarily switching a cell to Text format so numeric-looking strings
# are not auto-converted into numbers by Excel.
$defaultStyle = $ws.Cells.Item(1,1).Style

$ws.Cells.Item(2,4).Value2 = '43.132.23'
$ws.Cells.Item(2,5).Value2 = '  -1.42%  '

$ws.Cells.Item(3,4).Value2 = '2.272.08'
$ws.Cells.Item(3,5).Value2 = '  -0.84%  '

$ws.Cells.Item(4,5).Value2 = '  -0.24%  '

$c = $ws.Cells.Item(5,4)
$c.NumberFormat = "@"
$c.Value2 = '110.90'
$c.Style = $defaultStyle
$ws.Cells.Item(5,5).Value2 = '  -2.89%  '

$c = $ws.Cells.Item(6,4)
$c.NumberFormat = "@"
$c.Value2 = '263.51'
$c.Style = $defaultStyle
$ws.Cells.Item(6,5).Value2 = '  -1.84%  '

$c = $ws.Cells.Item(7,4)
$c.NumberFormat = "@"
$c.Value2 = '0.642'
$c.Style = $defaultStyle
$ws.Cells.Item(7,5).Value2 = '  +2.87%  '

$ws.Cells.Item(8,5).Value2 = '  -0.04%  '

$c = $ws.Cells.Item(9,4)
$c.NumberFormat = "@"
$c.Value2 = '0.605'
$c.Style = $defaultStyle
$ws.Cells.Item(9,5).Value2 = '  -2.24%  '

$c = $ws.Cells.Item(10,4)
$c.NumberFormat = "@"
$c.Value2 = '46.80'
$c.Style = $defaultStyle
$ws.Cells.Item(10,5).Value2 = '  -2.92%  '

$ws.Cells.Item(11,5).Value2 = '  -1.25%  '

$c = $ws.Cells.Item(12,4)
$c.NumberFormat = "@"
$c.Value2 = '9.09'
$c.Style = $defaultStyle
$ws.Cells.Item(12,5).Value2 = '  +2.97%  '

$ws.Cells.Item(13,5).Value2 = '  +1.64%  '

$c = $ws.Cells.Item(14,4)
$c.NumberFormat = "@"
$c.Value2 = '15.34'
$c.Style = $defaultStyle
$ws.Cells.Item(14,5).Value2 = '  -1.96%  '

$ws.Cells.Item(15,4).Value2 = '2.616.43'
$ws.Cells.Item(15,5).Value2 = '  -0.78%  '

$c = $ws.Cells.Item(16,4)
$c.NumberFormat = "@"
$c.Value2 = '0.860'
$c.Style = $defaultStyle
$ws.Cells.Item(16,5).Value2 = '  +1.66%  '

$ws.Cells.Item(17,4).Value2 = '2.274.00'
$ws.Cells.Item(17,5).Value2 = '  -0.78%  '

$ws.Cells.Item(18,4).Value2 = '43.103.68'
$ws.Cells.Item(18,5).Value2 = '  -1.28%  '

$ws.Cells.Item(19,5).Value2 = '  -1.78%  '

$c = $ws.Cells.Item(20,4)
$c.NumberFormat = "@"
$c.Value2 = '6.73'
$c.Style = $defaultStyle
$ws.Cells.Item(20,5).Value2 = '  +3.03%  '

$c = $ws.Cells.Item(21,4)
$c.NumberFormat = "@"
$c.Value2 = '71.73'
$c.Style = $defaultStyle
$ws.Cells.Item(21,5).Value2 = '  -0.93%  '

$c = $ws.Cells.Item(22,4)
$c.NumberFormat = "@"
$c.Value2 = '2.43'
$c.Style = $defaultStyle
$ws.Cells.Item(22,5).Value2 = '  -1.92%  '

$c = $ws.Cells.Item(23,4)
$c.NumberFormat = "@"
$c.Value2 = '233.73'
$c.Style = $defaultStyle

$c = $ws.Cells.Item(24,4)
$c.NumberFormat = "@"
$c.Value2 = '9.39'
$c.Style = $defaultStyle
$ws.Cells.Item(24,5).Value2 = '  -3.29%  '

$ws.Cells.Item(25,5).Value2 = '  +0.14%  '

$ws.Cells.Item(26,5).Value2 = '  +1.93%  '

$c = $ws.Cells.Item(27,4)
$c.NumberFormat = "@"
$c.Value2 = '11.29'
$c.Style = $defaultStyle
$ws.Cells.Item(27,5).Value2 = '  -3.03%  '

$c = $ws.Cells.Item(28,4)
$c.NumberFormat = "@"
$c.Value2 = '40.73'
$c.Style = $defaultStyle
$ws.Cells.Item(28,5).Value2 = '  -3.53%  '

$c = $ws.Cells.Item(29,4)
$c.NumberFormat = "@"
$c.Value2 = '3.35'
$c.Style = $defaultStyle
$ws.Cells.Item(29,5).Value2 = '  -1.65%  '

$c = $ws.Cells.Item(30,4)
$c.NumberFormat = "@"
$c.Value2 = '2.25'
$c.Style = $defaultStyle
$ws.Cells.Item(30,5).Value2 = '  -0.26%  '

$c = $ws.Cells.Item(31,4)
$c.NumberFormat = "@"
$c.Value2 = '173.31'
$c.Style = $defaultStyle
$ws.Cells.Item(31,5).Value2 = '  -2.04%  '

$c = $ws.Cells.Item(32,4)
$c.NumberFormat = "@"
$c.Value2 = '21.39'
$c.Style = $defaultStyle
$ws.Cells.Item(32,5).Value2 = '  -0.87%  '

$c = $ws.Cells.Item(33,4)
$c.NumberFormat = "@"
$c.Value2 = '0.0897'
$c.Style = $defaultStyle
$ws.Cells.Item(33,5).Value2 = '  -3.32%  '

$c = $ws.Cells.Item(34,4)
$c.NumberFormat = "@"
$c.Value2 = '5.64'
$c.Style = $defaultStyle
$ws.Cells.Item(34,5).Value2 = '  +1.17%  '

$c = $ws.Cells.Item(35,4)
$c.NumberFormat = "@"
$c.Value2 = '0.130'
$c.Style = $defaultStyle
$ws.Cells.Item(35,5).Value2 = '  +2.41%  '

$ws.Cells.Item(36,5).Value2 = '  +2.65%  '

$ws.Cells.Item(37,5).Value2 = '  -2.74%  '

$c = $ws.Cells.Item(38,4)
$c.NumberFormat = "@"
$c.Value2 = '3.91'
$c.Style = $defaultStyle
$ws.Cells.Item(38,5).Value2 = '  +3.21%  '

$ws.Cells.Item(39,5).Value2 = '  -3.68%  '

$c = $ws.Cells.Item(40,4)
$c.NumberFormat = "@"
$c.Value2 = '2.58'
$c.Style = $defaultStyle
$ws.Cells.Item(40,5).Value2 = '  +6.67%  '

$c = $ws.Cells.Item(41,4)
$c.NumberFormat = "@"
$c.Value2 = '14.25'
$c.Style = $defaultStyle
$ws.Cells.Item(41,5).Value2 = '  +1.87%  '

$c = $ws.Cells.Item(42,4)
$c.NumberFormat = "@"
$c.Value2 = '75.55'
$c.Style = $defaultStyle
$ws.Cells.Item(42,5).Value2 = '  +5.79%  '

$c = $ws.Cells.Item(43,4)
$c.NumberFormat = "@"
$c.Value2 = '0.236'
$c.Style = $defaultStyle
$ws.Cells.Item(43,5).Value2 = '  -3.03%  '

$c = $ws.Cells.Item(44,4)
$c.NumberFormat = "@"
$c.Value2 = '6.06'
$c.Style = $defaultStyle
$ws.Cells.Item(44,5).Value2 = '  -0.79%  '

$ws.Cells.Item(45,5).Value2 = '  -0.27%  '

$c = $ws.Cells.Item(46,4)
$c.NumberFormat = "@"
$c.Value2 = '1.36'
$c.Style = $defaultStyle
$ws.Cells.Item(46,5).Value2 = '  -4.26%  '

$c = $ws.Cells.Item(47,4)
$c.NumberFormat = "@"
$c.Value2 = '8.50'
$c.Style = $defaultStyle
$ws.Cells.Item(47,5).Value2 = '  -2.98%  '

$c = $ws.Cells.Item(48,4)
$c.NumberFormat = "@"
$c.Value2 = '1.25'
$c.Style = $defaultStyle
$ws.Cells.Item(48,5).Value2 = '  +2.05%  '

$c = $ws.Cells.Item(49,4)
$c.NumberFormat = "@"
$c.Value2 = '0.0995'
$c.Style = $defaultStyle
$ws.Cells.Item(49,5).Value2 = '  -1.28%  '

$c = $ws.Cells.Item(50,4)
$c.NumberFormat = "@"
$c.Value2 = '100.90'
$c.Style = $defaultStyle
$ws.Cells.Item(50,5).Value2 = '  -1.90%  '

$c = $ws.Cells.Item(51,4)
$c.NumberFormat = "@"
$c.Value2 = '0.600'
$c.Style = $defaultStyle
$ws.Cells.Item(51,5).Value2 = '  +9.56%  '
